$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the task description in A3 to include "and time"
$ws.Range("A3").Value = "Sort by location and time"

# Move the active selection to A7, matching the author's final cursor position
$ws.Range("A7").Select()
